$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Stamp the row formatting (style "1": centered, wrapped, bordered)
#    for every new row by copying the format of the last existing
#    data row (180), which already carries that exact style.
# ------------------------------------------------------------------
$ws.Range("A180:K180").Copy() | Out-Null
$newRows = 181..208
foreach ($r in $newRows) {
    $ws.Range("A" + $r + ":K" + $r).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Write the brand-new text values FIRST, in the exact order they
#    were originally authored, so the workbooks shared-string table
#    is appended in the same sequence (index 222..250) as the source.
# ------------------------------------------------------------------
$ws.Range("E181").Value = "majonez"  # sharedString 222
$ws.Range("D181").Value = "JAJKA NA TWARDO i KANAPKA Z HUMMUSEM"  # sharedString 223
$ws.Range("E184").Value = "hummus klasyczny"  # sharedString 224
$ws.Range("E185").Value = "szpinak"  # sharedString 225
$ws.Range("H185").Value = "1 garść"  # sharedString 226
$ws.Range("I181").Value = "zagotuj osoloną wodę; wrzuć jajka i nie przerywając wrzenia gotuj przez 8 min; ostudź i usuń skorupkę; na wierzch wyłóż majonez; chleb posmaruj hummusem; na wierzch wyłóż szpinak i szynkę"  # sharedString 227
$ws.Range("D188").Value = "JOGURT Z MUSLI I BANANEM"  # sharedString 228
$ws.Range("E189").Value = "musli z rodzynkami i orzechami"  # sharedString 229
$ws.Range("H189").Value = "8 łyżek"  # sharedString 230
$ws.Range("H190").Value = "9 łyżek"  # sharedString 231
$ws.Range("I188").Value = "wymieszaj jogurt z musli; dodaj pokrojonego banana"  # sharedString 232
$ws.Range("D192").Value = "MAKARON Z PESTO, SZPINAKIEM I KURCZAKIEM"  # sharedString 233
$ws.Range("E193").Value = "pesto zielone z bazylii"  # sharedString 234
$ws.Range("E194").Value = "makaron pełnoziarnisty"  # sharedString 235
$ws.Range("H194").Value = "1, 8 szklanki"  # sharedString 236
$ws.Range("E195").Value = "pierś z kurczaka"  # sharedString 237
$ws.Range("H195").Value = "2 porcje"  # sharedString 238
$ws.Range("E196").Value = "pomidory koktajlowe"  # sharedString 239
$ws.Range("E197").Value = "dynia, pestki, łuskane"  # sharedString 240
$ws.Range("H198").Value = "0.5 łyzki"  # sharedString 241
$ws.Range("I192").Value = "ugotuj makaron; mięso pokrój w kostkę, przypraw ziołami oraz solą; natrzyj olejem, zwiń w folię do pieczenia; upiecz w piekarniku (około 20 minut); makaron wymieszaj z pesto, szpinakiem, mięsem, pestkami dyni oraz pomidorkami koktajlowymi"  # sharedString 242
$ws.Range("D200").Value = "MORELE, SUSZONE"  # sharedString 243
$ws.Range("E200").Value = "morele suszone"  # sharedString 244
$ws.Range("H200").Value = "13 szt."  # sharedString 245
$ws.Range("D202").Value = "TWAROŻEK ZE SZCZYPIORKIEM"  # sharedString 246
$ws.Range("E204").Value = "ser twaróg chudy"  # sharedString 247
$ws.Range("H204").Value = "0,5 opakowania"  # sharedString 248
$ws.Range("H207").Value = "1 sztuka"  # sharedString 249
$ws.Range("I202").Value = "rozdrobnij widelcem twaróg i wymieszaj z twarogiem; dodaj posiekany szczypiorek; dopraw serek za pomocą ziół i pieprzu; podawaj z pieczywem z masłem"  # sharedString 250

# ------------------------------------------------------------------
# 3) Fill in the remaining cells for each new row: numeric fields,
#    and text fields that repeat a value used elsewhere.
# ------------------------------------------------------------------
# Row 181
$ws.Range("A181").Value = 35
$ws.Range("B181").Value = "śniadanie"
$ws.Range("C181").Value = 20
$ws.Range("F181").Value = 30
$ws.Range("G181").Value = "g"
$ws.Range("H181").Value = "1 łyżeczka"
$ws.Range("J181").Value = 606
$ws.Range("K181").Value = "Łukasz Wiertel"

# Row 182
$ws.Range("E182").Value = "jajko"
$ws.Range("F182").Value = 112
$ws.Range("G182").Value = "g"
$ws.Range("H182").Value = "2 szt."

# Row 183
$ws.Range("E183").Value = "chleb żytni razowy"
$ws.Range("F183").Value = 60
$ws.Range("G183").Value = "g"
$ws.Range("H183").Value = "2 kromki"

# Row 184
$ws.Range("F184").Value = 20
$ws.Range("G184").Value = "g"
$ws.Range("H184").Value = "2 łyżeczki"

# Row 185
$ws.Range("F185").Value = 25
$ws.Range("G185").Value = "g"

# Row 186
$ws.Range("E186").Value = "szynka z piersi kurczaka"
$ws.Range("F186").Value = 30
$ws.Range("G186").Value = "g"
$ws.Range("H186").Value = "2 plasterki"

# Row 187
$ws.Range("E187").Value = "X"

# Row 188
$ws.Range("A188").Value = 36
$ws.Range("B188").Value = "lunch"
$ws.Range("C188").Value = 5
$ws.Range("E188").Value = "banan"
$ws.Range("F188").Value = 120
$ws.Range("G188").Value = "g"
$ws.Range("H188").Value = "1 szt."
$ws.Range("J188").Value = 538
$ws.Range("K188").Value = "Łukasz Wiertel"

# Row 189
$ws.Range("F189").Value = 80
$ws.Range("G189").Value = "g"

# Row 190
$ws.Range("E190").Value = "jogurt naturalny"
$ws.Range("F190").Value = 180
$ws.Range("G190").Value = "g"

# Row 191
$ws.Range("E191").Value = "X"

# Row 192
$ws.Range("A192").Value = 37
$ws.Range("B192").Value = "obiad"
$ws.Range("C192").Value = 30
$ws.Range("E192").Value = "szpinak"
$ws.Range("F192").Value = 25
$ws.Range("G192").Value = "g"
$ws.Range("H192").Value = "1 garść"
$ws.Range("J192").Value = 900
$ws.Range("K192").Value = "Łukasz Wiertel"

# Row 193
$ws.Range("F193").Value = 40
$ws.Range("G193").Value = "g"
$ws.Range("H193").Value = "2 łyżki"

# Row 194
$ws.Range("F194").Value = 40
$ws.Range("G194").Value = "g"

# Row 195
$ws.Range("F195").Value = 200
$ws.Range("G195").Value = "g"

# Row 196
$ws.Range("F196").Value = 60
$ws.Range("G196").Value = "g"
$ws.Range("H196").Value = "3 szt."

# Row 197
$ws.Range("F197").Value = 10
$ws.Range("G197").Value = "g"
$ws.Range("H197").Value = "1 łyzka"

# Row 198
$ws.Range("E198").Value = "olej rzepakowy"
$ws.Range("F198").Value = 5
$ws.Range("G198").Value = "g"

# Row 199
$ws.Range("E199").Value = "X"

# Row 200
$ws.Range("A200").Value = 38
$ws.Range("B200").Value = "przekąska"
$ws.Range("C200").Value = 1
$ws.Range("F200").Value = 100
$ws.Range("G200").Value = "g"
$ws.Range("I200").Value = "Zjesc ze smakiem"
$ws.Range("J200").Value = 300
$ws.Range("K200").Value = "Łukasz Wiertel"

# Row 201
$ws.Range("E201").Value = "X"

# Row 202
$ws.Range("A202").Value = 39
$ws.Range("B202").Value = "kolacja"
$ws.Range("C202").Value = 7
$ws.Range("E202").Value = "masło"
$ws.Range("F202").Value = 10
$ws.Range("G202").Value = "g"
$ws.Range("H202").Value = "2 łyżeczki"
$ws.Range("J202").Value = 590
$ws.Range("K202").Value = "Łukasz Wiertel"

# Row 203
$ws.Range("E203").Value = "chleb żytni razowy"
$ws.Range("F203").Value = 120
$ws.Range("G203").Value = "g"
$ws.Range("H203").Value = "4 kromki"

# Row 204
$ws.Range("F204").Value = 100
$ws.Range("G204").Value = "g"

# Row 205
$ws.Range("E205").Value = "jogurt naturalny"
$ws.Range("F205").Value = 80
$ws.Range("G205").Value = "g"
$ws.Range("H205").Value = "4 łyżki"

# Row 206
$ws.Range("E206").Value = "szczypiorek "
$ws.Range("F206").Value = 10
$ws.Range("G206").Value = "g"
$ws.Range("H206").Value = "2 łyżeczki"

# Row 207
$ws.Range("E207").Value = "pomarańcze"
$ws.Range("F207").Value = 200
$ws.Range("G207").Value = "g"

# Row 208
$ws.Range("E208").Value = "X"

# ------------------------------------------------------------------
# 4) Row heights (matches the per-row ht values in the target sheet)
# ------------------------------------------------------------------
$ws.Rows.Item(181).RowHeight = 43.8
$ws.Rows.Item(182).RowHeight = 15
$ws.Rows.Item(183).RowHeight = 15
$ws.Rows.Item(184).RowHeight = 15
$ws.Rows.Item(185).RowHeight = 15
$ws.Rows.Item(186).RowHeight = 15
$ws.Rows.Item(187).RowHeight = 15
$ws.Rows.Item(188).RowHeight = 29.4
$ws.Rows.Item(189).RowHeight = 29.4
$ws.Rows.Item(190).RowHeight = 15
$ws.Rows.Item(191).RowHeight = 15
$ws.Rows.Item(192).RowHeight = 43.8
$ws.Rows.Item(193).RowHeight = 15
$ws.Rows.Item(194).RowHeight = 15
$ws.Rows.Item(195).RowHeight = 15
$ws.Rows.Item(196).RowHeight = 15
$ws.Rows.Item(197).RowHeight = 15
$ws.Rows.Item(198).RowHeight = 15
$ws.Rows.Item(199).RowHeight = 15
$ws.Rows.Item(200).RowHeight = 29.4
$ws.Rows.Item(201).RowHeight = 15
$ws.Rows.Item(202).RowHeight = 29.4
$ws.Rows.Item(203).RowHeight = 15
$ws.Rows.Item(204).RowHeight = 29.4
$ws.Rows.Item(205).RowHeight = 15
$ws.Rows.Item(206).RowHeight = 15
$ws.Rows.Item(207).RowHeight = 15
$ws.Rows.Item(208).RowHeight = 15

# ------------------------------------------------------------------
# 5) Update the active selection / scroll position to match the
#    final view state recorded in the workbook.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 196
$ws.Range("K203").Select() | Out-Null